$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("mmWave(InBed)")
$ws.Range("A38").Value = "'2026-01-28"
$ws.Range("B38").Value = "18:40:10"
$ws.Range("C38").Value = "18:00"
$ws.Range("D38").Value = "Bedroom"
$ws.Range("E38").Value = "In Bed"
$ws.Range("F38").Value = "Occupied"
$ws.Range("A39").Value = "'2026-01-28"
$ws.Range("B39").Value = "18:40:11"
$ws.Range("C39").Value = "18:00"
$ws.Range("D39").Value = "Bedroom"
$ws.Range("E39").Value = "In Bed"
$ws.Range("F39").Value = "Occupied"
$ws.Range("A40").Value = "'2026-01-28"
$ws.Range("B40").Value = "18:40:12"
$ws.Range("C40").Value = "18:00"
$ws.Range("D40").Value = "Bedroom"
$ws.Range("E40").Value = "In Bed"
$ws.Range("F40").Value = "Occupied"
$ws.Range("A41").Value = "'2026-01-28"
$ws.Range("B41").Value = "18:40:14"
$ws.Range("C41").Value = "18:00"
$ws.Range("D41").Value = "Bedroom"
$ws.Range("E41").Value = "In Bed"
$ws.Range("F41").Value = "Occupied"
$ws.Range("A42").Value = "'2026-01-28"
$ws.Range("B42").Value = "18:40:17"
$ws.Range("C42").Value = "18:00"
$ws.Range("D42").Value = "Bedroom"
$ws.Range("E42").Value = "In Bed"
$ws.Range("F42").Value = "Occupied"
$ws.Range("A43").Value = "'2026-01-28"
$ws.Range("B43").Value = "18:40:33"
$ws.Range("C43").Value = "18:00"
$ws.Range("D43").Value = "Bedroom"
$ws.Range("E43").Value = "In Bed"
$ws.Range("F43").Value = "Occupied"
$ws.Range("A44").Value = "'2026-01-28"
$ws.Range("B44").Value = "18:40:35"
$ws.Range("C44").Value = "18:00"
$ws.Range("D44").Value = "Bedroom"
$ws.Range("E44").Value = "In Bed"
$ws.Range("F44").Value = "Occupied"
$ws.Range("A45").Value = "'2026-01-28"
$ws.Range("B45").Value = "18:40:36"
$ws.Range("C45").Value = "18:00"
$ws.Range("D45").Value = "Bedroom"
$ws.Range("E45").Value = "In Bed"
$ws.Range("F45").Value = "Occupied"
$ws.Range("A46").Value = "'2026-01-28"
$ws.Range("B46").Value = "18:40:38"
$ws.Range("C46").Value = "18:00"
$ws.Range("D46").Value = "Bedroom"
$ws.Range("E46").Value = "In Bed"
$ws.Range("F46").Value = "Occupied"

$ws = $wb.Worksheets.Item("mmWave(BR)")
$ws.Range("A38").Value = "'2026-01-28"
$ws.Range("B38").Value = "18:40:11"
$ws.Range("C38").Value = "18:00"
$ws.Range("D38").Value = "Bedroom"
$ws.Range("E38").Value = 2
$ws.Range("F38").Value = "Occupied"
$ws.Range("A39").Value = "'2026-01-28"
$ws.Range("B39").Value = "18:40:11"
$ws.Range("C39").Value = "18:00"
$ws.Range("D39").Value = "Bedroom"
$ws.Range("E39").Value = 47
$ws.Range("F39").Value = "Occupied"
$ws.Range("A40").Value = "'2026-01-28"
$ws.Range("B40").Value = "18:40:12"
$ws.Range("C40").Value = "18:00"
$ws.Range("D40").Value = "Bedroom"
$ws.Range("E40").Value = 41
$ws.Range("F40").Value = "Occupied"
$ws.Range("A41").Value = "'2026-01-28"
$ws.Range("B41").Value = "18:40:14"
$ws.Range("C41").Value = "18:00"
$ws.Range("D41").Value = "Bedroom"
$ws.Range("E41").Value = 2
$ws.Range("F41").Value = "Occupied"
$ws.Range("A42").Value = "'2026-01-28"
$ws.Range("B42").Value = "18:40:17"
$ws.Range("C42").Value = "18:00"
$ws.Range("D42").Value = "Bedroom"
$ws.Range("E42").Value = 1
$ws.Range("F42").Value = "Occupied"
$ws.Range("A43").Value = "'2026-01-28"
$ws.Range("B43").Value = "18:40:34"
$ws.Range("C43").Value = "18:00"
$ws.Range("D43").Value = "Bedroom"
$ws.Range("E43").Value = 40
$ws.Range("F43").Value = "Occupied"
$ws.Range("A44").Value = "'2026-01-28"
$ws.Range("B44").Value = "18:40:35"
$ws.Range("C44").Value = "18:00"
$ws.Range("D44").Value = "Bedroom"
$ws.Range("E44").Value = 10
$ws.Range("F44").Value = "Occupied"
$ws.Range("A45").Value = "'2026-01-28"
$ws.Range("B45").Value = "18:40:36"
$ws.Range("C45").Value = "18:00"
$ws.Range("D45").Value = "Bedroom"
$ws.Range("E45").Value = 21
$ws.Range("F45").Value = "Occupied"
$ws.Range("A46").Value = "'2026-01-28"
$ws.Range("B46").Value = "18:40:38"
$ws.Range("C46").Value = "18:00"
$ws.Range("D46").Value = "Bedroom"
$ws.Range("E46").Value = 3
$ws.Range("F46").Value = "Occupied"

$ws = $wb.Worksheets.Item("mmWave(HR)")
$ws.Range("A38").Value = "'2026-01-28"
$ws.Range("B38").Value = "18:40:10"
$ws.Range("C38").Value = "18:00"
$ws.Range("D38").Value = "Bedroom"
$ws.Range("E38").Value = 50
$ws.Range("F38").Value = "Occupied"
$ws.Range("A39").Value = "'2026-01-28"
$ws.Range("B39").Value = "18:40:11"
$ws.Range("C39").Value = "18:00"
$ws.Range("D39").Value = "Bedroom"
$ws.Range("E39").Value = 95
$ws.Range("F39").Value = "Occupied"
$ws.Range("A40").Value = "'2026-01-28"
$ws.Range("B40").Value = "18:40:12"
$ws.Range("C40").Value = "18:00"
$ws.Range("D40").Value = "Bedroom"
$ws.Range("E40").Value = 89
$ws.Range("F40").Value = "Occupied"
$ws.Range("A41").Value = "'2026-01-28"
$ws.Range("B41").Value = "18:40:14"
$ws.Range("C41").Value = "18:00"
$ws.Range("D41").Value = "Bedroom"
$ws.Range("E41").Value = 50
$ws.Range("F41").Value = "Occupied"
$ws.Range("A42").Value = "'2026-01-28"
$ws.Range("B42").Value = "18:40:17"
$ws.Range("C42").Value = "18:00"
$ws.Range("D42").Value = "Bedroom"
$ws.Range("E42").Value = 49
$ws.Range("F42").Value = "Occupied"
$ws.Range("A43").Value = "'2026-01-28"
$ws.Range("B43").Value = "18:40:34"
$ws.Range("C43").Value = "18:00"
$ws.Range("D43").Value = "Bedroom"
$ws.Range("E43").Value = 88
$ws.Range("F43").Value = "Occupied"
$ws.Range("A44").Value = "'2026-01-28"
$ws.Range("B44").Value = "18:40:35"
$ws.Range("C44").Value = "18:00"
$ws.Range("D44").Value = "Bedroom"
$ws.Range("E44").Value = 58
$ws.Range("F44").Value = "Occupied"
$ws.Range("A45").Value = "'2026-01-28"
$ws.Range("B45").Value = "18:40:36"
$ws.Range("C45").Value = "18:00"
$ws.Range("D45").Value = "Bedroom"
$ws.Range("E45").Value = 69
$ws.Range("F45").Value = "Occupied"
$ws.Range("A46").Value = "'2026-01-28"
$ws.Range("B46").Value = "18:40:38"
$ws.Range("C46").Value = "18:00"
$ws.Range("D46").Value = "Bedroom"
$ws.Range("E46").Value = 51
$ws.Range("F46").Value = "Occupied"
